$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.426.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.189.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.17%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.189.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.511"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.714.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.511.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.186.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.112"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  +6.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "510.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0899"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0425"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.126"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.12%  "
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("E43").Value = "  -4.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0675"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.854.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.09%  "
